$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44245
$ws.Range("J2").Value = 9000
$ws.Range("K2").Value = 3000
$ws.Range("M2").Value = 3000
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 30

# Row 3
$ws.Range("D3").Value = 44245
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2500
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 25

# Row 4
$ws.Range("D4").Value = 44231
$ws.Range("J4").Value = 12000

# Row 5
$ws.Range("D5").Value = 44874
$ws.Range("J5").Value = 7900

# Row 6
$ws.Range("D6").Value = 44875
$ws.Range("J6").Value = 7900

# Row 7
$ws.Range("D7").Value = 44162
$ws.Range("J7").Value = 7000
$ws.Range("O7").Value = "Provincia de Chacabuco"

# Row 8
$ws.Range("D8").Value = 44210
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 8800
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 2750
$ws.Range("O8").Value = "Provincia de Chacabuco"
$ws.Range("P8").Value = 28

# Row 9
$ws.Range("D9").Value = 44229
$ws.Range("J9").Value = 16000

# Row 10
$ws.Range("D10").Value = 44232
$ws.Range("J10").Value = 16000

# Row 11
$ws.Range("D11").Value = 44902
$ws.Range("J11").Value = 7000

# Row 12
$ws.Range("D12").Value = 44602
$ws.Range("J12").Value = 12000

# Row 13
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 6000
$ws.Range("K13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2500
$ws.Range("P13").Value = 25

# Row 14
$ws.Range("D14").Value = 44159
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 7000
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = 3000
$ws.Range("P14").Value = 30

# Row 15
$ws.Range("D15").Value = 44876
$ws.Range("J15").Value = 7900

# Row 16
$ws.Range("D16").Value = 44168
$ws.Range("J16").Value = 7000

# Row 17
$ws.Range("D17").Value = 44215
$ws.Range("J17").Value = 16000

# Row 18
$ws.Range("D18").Value = 44204

# Row 19
$ws.Range("D19").Value = 44189
$ws.Range("J19").Value = 16000

# Row 20
$ws.Range("D20").Value = 44859

# Row 21
$ws.Range("D21").Value = 44873

# Row 22
$ws.Range("D22").Value = 44209
$ws.Range("J22").Value = 7000
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = 2750
$ws.Range("P22").Value = 28

# Row 23
$ws.Range("D23").Value = 44880

# Row 24
$ws.Range("D24").Value = 44230
$ws.Range("J24").Value = 16000

# Row 25
$ws.Range("D25").Value = 44881
$ws.Range("J25").Value = 7900

# Row 26
$ws.Range("D26").Value = 44188
$ws.Range("J26").Value = 12000

# Row 27
$ws.Range("D27").Value = 44181
$ws.Range("J27").Value = 12000

# Row 28
$ws.Range("D28").Value = 44187
$ws.Range("J28").Value = 12000

# Row 29
$ws.Range("D29").Value = 44883
$ws.Range("J29").Value = 9700

# Row 30
$ws.Range("D30").Value = 44160
$ws.Range("J30").Value = 7000

# Row 31
$ws.Range("D31").Value = 44166
$ws.Range("J31").Value = 7000

# Row 32
$ws.Range("D32").Value = 44860
$ws.Range("J32").Value = 7900

# Row 33
$ws.Range("D33").Value = 44882
$ws.Range("J33").Value = 7900

# Row 34
$ws.Range("D34").Value = 44186
$ws.Range("J34").Value = 10000

# Row 35
$ws.Range("D35").Value = 44167
$ws.Range("J35").Value = 7000
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = 3000
$ws.Range("O35").Value = "Provincia de Chacabuco"
$ws.Range("P35").Value = 30

# Row 36
$ws.Range("D36").Value = 44847
$ws.Range("J36").Value = 7900
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = 3000
$ws.Range("P36").Value = 30

# Row 37
$ws.Range("D37").Value = 44846
$ws.Range("J37").Value = 7900

# Row 38
$ws.Range("D38").Value = 44901
$ws.Range("J38").Value = 7000

# Row 39
$ws.Range("D39").Value = 44855
$ws.Range("J39").Value = 7900

# Row 40
$ws.Range("D40").Value = 44161
$ws.Range("J40").Value = 7000

# Row 41
$ws.Range("D41").Value = 44600
$ws.Range("J41").Value = 1300
$ws.Range("K41").Value = 3500
$ws.Range("L41").Value = 4000
$ws.Range("M41").Value = 3808
$ws.Range("O41").Value = "Región Metropolitana"
$ws.Range("P41").Value = 38

# Row 42
$ws.Range("D42").Value = 44845
$ws.Range("J42").Value = 7900

# Row 43
$ws.Range("D43").Value = 44214
$ws.Range("J43").Value = 7000
